$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits at the end of
#    the "You need the materials..." paragraph.
# ---------------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Rework the "What does clicking the Submit button do?" heading and the
#    paragraph that follows it (which used to be a Heading2 "<form>" block
#    with the _rfgpcdmwumw9 bookmark). The heading's paragraph-mark run
#    properties gain a color, the bookmark moves to the end of the heading
#    paragraph, and the following paragraph becomes a plain "This is our
#    code:" paragraph.
# ---------------------------------------------------------------------------
$headingPara = $null
$formPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "What does clicking the Submit button do?*") {
        $headingPara = $p
        $formPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

$rng = $d.Range($headingPara.Range.Start, $formPara.Range.End)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:color w:val="000000"/><w:sz w:val="30"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="30"/></w:rPr><w:t xml:space="preserve">What does </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="30"/></w:rPr><w:t>clicking</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="30"/></w:rPr><w:t xml:space="preserve"> the Submit button do?</w:t></w:r><w:bookmarkStart w:id="101" w:name="_rfgpcdmwumw9" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="101"/></w:p><w:p><w:pPr><w:rPr><w:sz w:val="30"/></w:rPr></w:pPr><w:r><w:t>This is our code:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) The final list-paragraph ("When running your html website...") loses
#    its ListParagraph/numPr formatting and instead gets the "_GoBack"
#    bookmark placed right at its start.
# ---------------------------------------------------------------------------
$lastPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*hen running your html website*") {
        $lastPara = $p
        break
    }
}

$rng3 = $lastPara.Range
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="201" w:name="_GoBack"/><w:bookmarkEnd w:id="201"/><w:r><w:t>W</w:t></w:r><w:r><w:t xml:space="preserve">hen running your html website, there will appear a button. Due to you not having the correct page yet, it will prompt with an error. But otherwise, it would use another file and record the text that you input in the second file.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng3.InsertXML($xml3)

Write-Host "All edits applied"
